# Weekly update: insert a new price-record row for the week of 2021-11-16
# (row 179) in the "Coliflor" sheet, pushing the existing rows 179-187 down
# to 180-188 (dimension grows from A1:R187 to A1:R188).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 179; formatting (e.g. the date
# number format on column D) is inherited from the row below, matching how
# the rest of the table is formatted.
$ws.Rows("179:179").Insert()

$ws.Cells.Item(179, 1).Value = 7
$ws.Cells.Item(179, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(179, 3).Value = "Ñuble"
$ws.Cells.Item(179, 4).Value = 44516
$ws.Cells.Item(179, 5).Value = 16
$ws.Cells.Item(179, 6).Value = 100112008
$ws.Cells.Item(179, 7).Value = "Coliflor"
$ws.Cells.Item(179, 8).Value = "Sin especificar"
$ws.Cells.Item(179, 9).Value = "Primera"
$ws.Cells.Item(179, 10).Value = 360
$ws.Cells.Item(179, 11).Value = 700
$ws.Cells.Item(179, 12).Value = 750
$ws.Cells.Item(179, 13).Value = 725
$ws.Cells.Item(179, 14).Value = "$/unidad"
$ws.Cells.Item(179, 15).Value = "Región del Maule"
$ws.Cells.Item(179, 16).Value = 725
$ws.Cells.Item(179, 17).Value = 1
$ws.Cells.Item(179, 18).Value = "Hortaliza"
